$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trans_SA_FAresults")

# Header row: insert new "LB" label in D1 (pushing the conceptual labels so
# obj/CI/train_time/test_time line up with columns E:H), and add the
# previously-missing "test_time" label in H1.
$ws.Range("E1").Value = "obj"
$ws.Range("F1").Value = "CI"
$ws.Range("G1").Value = "train_time"
$ws.Range("H1").Value = "test_time"
$ws.Range("D1").Value = "LB"

# Updated numeric results for columns D:H, rows 2:11 (values only change in place,
# no columns/rows are actually inserted or moved).
$data = @(
    @(2381917.3708052002, 2224109.5738952402, 128755.877329667,   172.79141712188701, 7.2834231853485099),
    @(2821117.4904798502, 2702340.9756223098, 158119.72047798801, 134.83133006095801, 6.7389283180236799),
    @(3124434.31244912,   2985251.8338046698, 154233.01989820701, 318.20528793334898, 10.941272974014201),
    @(3208616.19458741,   3119646.2051250502, 195906.26697943901, 182.00984430313099, 7.4736657142639098),
    @(3279854.349891,     3100418.8371999399, 183713.694034418,   99.673047065734806, 6.2072639465331996),
    @(2526503.0701657198, 2446379.7740750802, 214354.718946526,   76.719655990600501, 7.2813510894775302),
    @(2829311.4969536001, 2765838.6373697198, 195737.92233664601, 66.538102626800494, 6.9251949787139804),
    @(3009202.9975562799, 3025749.1541551198, 259595.53274857899, 87.090435028076101, 7.5620210170745796),
    @(3081331.1331428899, 3083036.93960738,   268750.53732023702, 81.368405103683401, 7.57218289375305),
    @(3124531.7976661902, 3078117.6668686401, 263259.74107469001, 61.789297342300401, 6.77239966392517)
)

$row = 2
foreach ($r in $data) {
    $ws.Range("D$row").Value = $r[0]
    $ws.Range("E$row").Value = $r[1]
    $ws.Range("F$row").Value = $r[2]
    $ws.Range("G$row").Value = $r[3]
    $ws.Range("H$row").Value = $r[4]
    $row++
}

# Keep the selection consistent with the authored file.
$null = $ws.Range("D1").Select()
